# EWN: More changes after first full play
# Updates item-rule text on the "Cards" sheet and sticker/rulebook text on
# the "Stickers" sheet, then restores the active-cell selections.

$wb = $excel.ActiveWorkbook
$cards = $wb.Worksheets.Item("Cards")
$stickers = $wb.Worksheets.Item("Stickers")

# ---- Cards sheet: item rule text tweaks ----

# Bridge's item rule gained "Return them." at the end (appears 3 times:
# paired with Broom, Biscuit, and Crowbar).
$cards.Range("I10").Value = "Look at the top 2 cards of the deck. Activate an item on one of them. Return them."
$cards.Range("I12").Value = "Look at the top 2 cards of the deck. Activate an item on one of them. Return them."
$cards.Range("I13").Value = "Look at the top 2 cards of the deck. Activate an item on one of them. Return them."

# Stick Insect's item rule was replaced outright (appears paired with
# Daisies, Butterfly, and Fig Leaf).
$cards.Range("I15").Value = "Discard a card. Draw 2 cards. Discard a card."
$cards.Range("I17").Value = "Discard a card. Draw 2 cards. Discard a card."
$cards.Range("F19").Value = "Discard a card. Draw 2 cards. Discard a card."

# Fig Leaf's item rule reworded (appears paired with Daisies, Butterfly,
# and Stick Insect).
$cards.Range("I16").Value = "Draw 2 cards. Return a card. Return one of your Equipment cards."
$cards.Range("I18").Value = "Draw 2 cards. Return a card. Return one of your Equipment cards."
$cards.Range("I19").Value = "Draw 2 cards. Return a card. Return one of your Equipment cards."

# ---- Stickers sheet: sticker costs / text tweaks ----

# Sticker s1 (row 2)
$stickers.Range("F2").Value = "3{+}:\nDraw a card. Return or discard a card. Equip a card and activate it."
$stickers.Range("G2").Value = "3{+}:\nScore 1 of your Equipment cards as a fragment. Draw a card"
$stickers.Range("H2").Value = "3{+}:\nAttach this as a tab to the edge of the card. While it is in the deck, at any time you may discard every card above this."
$stickers.Range("I2").Value = "2{+}:\nEither steal a fading card, or score a card normally."
$stickers.Range("J2").Value = "3{+}:\nAlways: When scoring treat this as any item.\nActivate: Rotate any number of equipment cards"
$stickers.Range("L2").Value = "A: Play until you score every card from Era B and get at least 18 points during a single game, then pull the second sticker from the sleeve and apply its rules to the rule document."
$stickers.Range("M2").Value = "B: After a game you may spend from a card the number of {+} on a mod sticker to apply it to the card (no covering other stickers). You may also spend {+} to mark the back of a card with 1 letter."

# Sticker s2 (row 3)
$stickers.Range("G3").Value = "2{+}:\n Return 2 cards. Draw 4 cards. Discard or return 2 cards."
$stickers.Range("H3").Value = "3{+}:\nRotate, Score, or Discard and repace the paradox."
$stickers.Range("J3").Value = "3{+}:\n Look at up to 5 cards from the top of the deck. Return them in any order."
$stickers.Range("L3").Value = "A: Play until you score 16 points from cards in Era C during a single game, then pull the final sticker from the sleeve and apply its rules to the rule document."
$stickers.Range("M3").Value = "C: Apply the Time Rift sticker to the card with the most {-} on it. Add 4{-} to the card. You may spend {+}{-} to change the name of an item on a card (which changes how it is scored)."

# Sticker s3 (row 4)
$stickers.Range("F4").Value = "5{+}:\nLook at the top 4 cards of the deck. Score or discard any number of them. Return the remainder in any order."
$stickers.Range("I4").Value = "3{+}:\n Draw 4 cards, return or discard 4 cards. Rotate any 1 Equipment card. "
$stickers.Range("J4").Value = "6{+}:\n Always: you may activate 3 items a turn instead of 2."
$stickers.Range("L4").Value = "A: Apply the Singularity sticker to the other half of the card with the Time Rift on it. Add 4{-} to the card. Play until you score both the Singularity and 30+ points in a single game. Then you win."
$stickers.Range("M4").Value = "D:You may spend {+} to +1 a number on a card, or {+}{-}{-} to -1 a number on a card. You may spend {+}{+} to mark on the back of a card where it starts the game (hand, paradox, num in deck, etc.)."

# ---- Restore active-cell selections (Cards first, Stickers last so the
#      workbook ends with Stickers as the active/tab-selected sheet, as
#      it was originally). ----
$cards.Range("I4").Select()
$stickers.Range("I5").Select()
